$wb = $excel.ActiveWorkbook

# --- targetAssoc sheet: rebuild the table with the new headers/rows ---
$ws = $wb.Worksheets.Item("targetAssoc")

# Column E (diseaseCount) gets its own best-fit width now that it's the last
# column; set it on the about-to-be-deleted column F so the width carries
# over to E once the F/G columns shift out.
$ws.Columns.Item(6).ColumnWidth = 15.6

# Drop the old F (outputFilename) / G (' ') columns entirely.
$ws.Columns("F:G").Delete()

# New header row (renamed columns).
$ws.Range("A1").Value = "suffix_Url"
$ws.Range("B1").Value = "target_ID"
$ws.Range("C1").Value = "target_Name"
$ws.Range("D1").Value = "PMTLcode"
$ws.Range("E1").Value = "diseaseCount"

# Row 2 stays ALK, but drop the old F2 "outputFilename" value (done above via column delete).
$ws.Range("A2").Value = "/target/ENSG00000171094/associations"
$ws.Range("B2").Value = "ENSG00000171094"
$ws.Range("C2").Value = "ALK"
$ws.Range("D2").Value = "Relevant Molecular Target"
$ws.Range("E2").Value = 758

# New rows 3-5, matching the style already used on row 2 (column A uses the
# "s=2" style, columns B-E use the "s=1" style) by copying that formatting down.
$ws.Range("A2:E2").Copy($ws.Range("A3:E3"))
$ws.Range("A2:E2").Copy($ws.Range("A4:E4"))
$ws.Range("A2:E2").Copy($ws.Range("A5:E5"))
$excel.CutCopyMode = 0
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75

$ws.Range("A3").Value = "/target/ENSG00000133703/associations"
$ws.Range("B3").Value = "ENSG00000133703"
$ws.Range("C3").Value = "KRAS"
$ws.Range("D3").Value = "Relevant Molecular Target"
$ws.Range("E3").Value = 0

$ws.Range("A4").Value = "/target/ENSG00000232810/associations"
$ws.Range("B4").Value = "ENSG00000232810"
$ws.Range("C4").Value = "TNF"
$ws.Range("D4").Value = "Unspecified Target"
$ws.Range("E4").Value = 2676

$ws.Range("A5").Value = "/target/ENSG00000169083/associations"
$ws.Range("B5").Value = "ENSG00000169083"
$ws.Range("C5").Value = "AR"
$ws.Range("D5").Value = "Non-Relevant Molecular Target"
$ws.Range("E5").Value = 1022

# Touch the font table so the "family" attribute gets populated (matches the
# re-saved workbook's style sheet).
$ws.Range("B1").Font.Family = 2

$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null

# --- copy sheet: selection only covers through column E now ---
$ws2 = $wb.Worksheets.Item("copy")
$ws2.Activate() | Out-Null
$ws2.Range("A2:E5").Select() | Out-Null

$ws.Activate() | Out-Null
